$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 67
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 3).NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2025-10-09"
$ws.Cells.Item($row, 2).Value = "15:32:09"
$ws.Cells.Item($row, 3).Value = "1.00 EUR = 1,776.4528"
